$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on cells whose new values would otherwise be
# auto-detected as numeric by Excel (single-dot decimal-looking strings),
# so they remain text exactly as in the source data.
$textCells = @("D5","D6","D7","D10","D11","D16","D19","D21","D22","D23","D25","D28","D31","D32","D34","D36","D37","D38","D39","D40","D41","D42","D43","D46","D49","D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "42.593.55"
$ws.Range("E2").Value = "  -7.68%  "
$ws.Range("D3").Value = "2.551.51"
$ws.Range("E3").Value = "  -1.84%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "295.83"
$ws.Range("E5").Value = "  -5.04%  "
$ws.Range("D6").Value = "91.03"
$ws.Range("E6").Value = "  -8.02%  "
$ws.Range("D7").Value = "0.572"
$ws.Range("E7").Value = "  -4.47%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -5.92%  "
$ws.Range("D10").Value = "35.45"
$ws.Range("E10").Value = "  -8.94%  "
$ws.Range("D11").Value = "0.0804"
$ws.Range("E11").Value = "  -4.23%  "
$ws.Range("E12").Value = "  -6.31%  "
$ws.Range("D13").Value = "2.940.27"
$ws.Range("E13").Value = "  -1.89%  "
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "2.540.17"
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("D16").Value = "0.861"
$ws.Range("E16").Value = "  -5.88%  "
$ws.Range("E17").Value = "  -5.28%  "
$ws.Range("D18").Value = "42.586.64"
$ws.Range("E18").Value = "  -7.74%  "
$ws.Range("D19").Value = "6.66"
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("D20").Value = "0.0₃0966"
$ws.Range("E20").Value = "  -5.16%  "
$ws.Range("D21").Value = "12.49"
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("D22").Value = "72.41"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "257.80"
$ws.Range("E23").Value = "  -10.32%  "
$ws.Range("E24").Value = "  -5.80%  "
$ws.Range("D25").Value = "29.48"
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("E26").Value = "  -5.76%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "9.96"
$ws.Range("E28").Value = "  -7.12%  "
$ws.Range("E29").Value = "  -4.42%  "
$ws.Range("E30").Value = "  -5.15%  "
$ws.Range("D31").Value = "5.88"
$ws.Range("E31").Value = "  -5.23%  "
$ws.Range("D32").Value = "150.59"
$ws.Range("E32").Value = "  -3.15%  "
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("D34").Value = "3.40"
$ws.Range("E34").Value = "  -5.27%  "
$ws.Range("E35").Value = "  -3.01%  "
$ws.Range("D36").Value = "0.0789"
$ws.Range("E36").Value = "  -5.89%  "
$ws.Range("D37").Value = "0.114"
$ws.Range("E37").Value = "  -7.36%  "
$ws.Range("D38").Value = "24.21"
$ws.Range("E38").Value = "  +6.43%  "
$ws.Range("D39").Value = "0.119"
$ws.Range("E39").Value = "  -3.67%  "
$ws.Range("D40").Value = "15.82"
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("D41").Value = "3.42"
$ws.Range("E41").Value = "  -4.46%  "
$ws.Range("D42").Value = "0.0308"
$ws.Range("E42").Value = "  -6.60%  "
$ws.Range("D43").Value = "3.80"
$ws.Range("E43").Value = "  -4.11%  "
$ws.Range("D44").Value = "2.068.33"
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Value = "84.27"
$ws.Range("E46").Value = "  -13.11%  "
$ws.Range("E47").Value = "  +2.83%  "
$ws.Range("D48").Value = "2.794.72"
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("D49").Value = "8.72"
$ws.Range("E49").Value = "  -9.53%  "
$ws.Range("E50").Value = "  -3.66%  "
$ws.Range("D51").Value = "103.12"
$ws.Range("E51").Value = "  -5.16%  "
